# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest values from the scheduled data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$price = $ws.Range("D2")
$price.NumberFormat = "@"
$price.Value = "72.326.56"
$price.ClearFormats()
$ws.Range("E2").Value = "  +4.37%  "
$price = $ws.Range("D3")
$price.NumberFormat = "@"
$price.Value = "3.618.43"
$price.ClearFormats()
$ws.Range("E3").Value = "  +6.70%  "
$price = $ws.Range("D4")
$price.NumberFormat = "@"
$price.Value = "1.00"
$price.ClearFormats()
$ws.Range("E4").Value = "  +0.15%  "
$price = $ws.Range("D5")
$price.NumberFormat = "@"
$price.Value = "598.95"
$price.ClearFormats()
$ws.Range("E5").Value = "  +1.89%  "
$price = $ws.Range("D6")
$price.NumberFormat = "@"
$price.Value = "183.25"
$price.ClearFormats()
$ws.Range("E6").Value = "  +1.61%  "
$price = $ws.Range("D7")
$price.NumberFormat = "@"
$price.Value = "3.608.44"
$price.ClearFormats()
$ws.Range("E7").Value = "  +6.71%  "
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +5.43%  "
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("E12").Value = "  +4.34%  "
$price = $ws.Range("D13")
$price.NumberFormat = "@"
$price.Value = "0.0000291"
$price.ClearFormats()
$ws.Range("E13").Value = "  +2.92%  "
$price = $ws.Range("D14")
$price.NumberFormat = "@"
$price.Value = "708.14"
$price.ClearFormats()
$ws.Range("E14").Value = "  +4.26%  "
$price = $ws.Range("D15")
$price.NumberFormat = "@"
$price.Value = "4.202.91"
$price.ClearFormats()
$ws.Range("E15").Value = "  +6.89%  "
$price = $ws.Range("D16")
$price.NumberFormat = "@"
$price.Value = "8.99"
$price.ClearFormats()
$ws.Range("E16").Value = "  +3.93%  "
$price = $ws.Range("D17")
$price.NumberFormat = "@"
$price.Value = "72.396.18"
$price.ClearFormats()
$ws.Range("E17").Value = "  +4.34%  "
$price = $ws.Range("D18")
$price.NumberFormat = "@"
$price.Value = "3.613.11"
$price.ClearFormats()
$ws.Range("E18").Value = "  +6.47%  "
$ws.Range("E19").Value = "  +1.84%  "
$price = $ws.Range("D20")
$price.NumberFormat = "@"
$price.Value = "18.62"
$price.ClearFormats()
$ws.Range("E20").Value = "  +4.92%  "
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("E22").Value = "  +3.47%  "
$price = $ws.Range("D23")
$price.NumberFormat = "@"
$price.Value = "5.90"
$price.ClearFormats()
$ws.Range("E23").Value = "  +8.63%  "
$price = $ws.Range("D24")
$price.NumberFormat = "@"
$price.Value = "17.91"
$price.ClearFormats()
$ws.Range("E24").Value = "  +4.42%  "
$price = $ws.Range("D25")
$price.NumberFormat = "@"
$price.Value = "105.30"
$price.ClearFormats()
$ws.Range("E25").Value = "  +2.21%  "
$price = $ws.Range("D26")
$price.NumberFormat = "@"
$price.Value = "4.05"
$price.ClearFormats()
$ws.Range("E26").Value = "  +2.91%  "
$price = $ws.Range("D27")
$price.NumberFormat = "@"
$price.Value = "2.85"
$price.ClearFormats()
$ws.Range("E27").Value = "  +4.45%  "
$price = $ws.Range("D28")
$price.NumberFormat = "@"
$price.Value = "9.99"
$price.ClearFormats()
$ws.Range("E28").Value = "  +3.88%  "
$price = $ws.Range("D29")
$price.NumberFormat = "@"
$price.Value = "35.54"
$price.ClearFormats()
$ws.Range("E29").Value = "  +4.65%  "
$price = $ws.Range("D30")
$price.NumberFormat = "@"
$price.Value = "9.19"
$price.ClearFormats()
$ws.Range("E30").Value = "  +4.98%  "
$price = $ws.Range("D31")
$price.NumberFormat = "@"
$price.Value = "7.46"
$price.ClearFormats()
$ws.Range("E31").Value = "  +7.45%  "
$price = $ws.Range("D32")
$price.NumberFormat = "@"
$price.Value = "4.19"
$price.ClearFormats()
$ws.Range("E32").Value = "  +16.37%  "
$price = $ws.Range("D33")
$price.NumberFormat = "@"
$price.Value = "592.47"
$price.ClearFormats()
$ws.Range("E33").Value = "  +6.36%  "
$price = $ws.Range("D34")
$price.NumberFormat = "@"
$price.Value = "11.37"
$price.ClearFormats()
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  +1.53%  "
$price = $ws.Range("D36")
$price.NumberFormat = "@"
$price.Value = "59.63"
$price.ClearFormats()
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +3.91%  "
$price = $ws.Range("D39")
$price.NumberFormat = "@"
$price.Value = "3.639.76"
$price.ClearFormats()
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  +8.84%  "
$price = $ws.Range("D41")
$price.NumberFormat = "@"
$price.Value = "36.08"
$price.ClearFormats()
$ws.Range("E41").Value = "  +0.87%  "
$price = $ws.Range("D42")
$price.NumberFormat = "@"
$price.Value = "3.49"
$price.ClearFormats()
$ws.Range("E42").Value = "  +6.82%  "
$price = $ws.Range("D43")
$price.NumberFormat = "@"
$price.Value = "2.81"
$price.ClearFormats()
$ws.Range("E43").Value = "  +4.47%  "
$ws.Range("E44").Value = "  +6.59%  "
$price = $ws.Range("D45")
$price.NumberFormat = "@"
$price.Value = "0.350"
$price.ClearFormats()
$ws.Range("E45").Value = "  +2.69%  "
$price = $ws.Range("D46")
$price.NumberFormat = "@"
$price.Value = "3.43"
$price.ClearFormats()
$ws.Range("E46").Value = "  +3.64%  "
$price = $ws.Range("D47")
$price.NumberFormat = "@"
$price.Value = "2.81"
$price.ClearFormats()
$ws.Range("E47").Value = "  +4.54%  "
$price = $ws.Range("D48")
$price.NumberFormat = "@"
$price.Value = "1.48"
$price.ClearFormats()
$ws.Range("E48").Value = "  +4.97%  "
$ws.Range("E49").Value = "  +1.89%  "
$price = $ws.Range("D50")
$price.NumberFormat = "@"
$price.Value = "0.999"
$price.ClearFormats()
$ws.Range("E50").Value = "  -0.25%  "
$price = $ws.Range("D51")
$price.NumberFormat = "@"
$price.Value = "133.82"
$price.ClearFormats()
$ws.Range("E51").Value = "  +0.19%  "
